$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '97.565.83'
$ws.Range("E2").Value = '  +4.17%  '

$ws.Range("D3").Value = '3.302.61'
$ws.Range("E3").Value = '  +6.78%  '

$ws.Range("D4").Value = '''1.00'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").Value = '''243.61'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.89%  '

$ws.Range("D6").Value = '''619.58'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.38%  '

$ws.Range("D7").Value = '''1.14'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +1.59%  '

$ws.Range("D8").Value = '''0.386'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.31%  '

$ws.Range("D9").Value = '''1.00'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.14%  '

$ws.Range("D10").Value = '3.294.91'
$ws.Range("E10").Value = '  +6.58%  '

$ws.Range("D11").Value = '''0.804'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.96%  '

$ws.Range("E12").Value = '  +1.61%  '

$ws.Range("D13").Value = '97.879.57'
$ws.Range("E13").Value = '  +4.88%  '

$ws.Range("D14").Value = '''0.0000245'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.14%  '

$ws.Range("D15").Value = '''35.22'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +2.07%  '

$ws.Range("B16").Value = 'Toncoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D16").Value = '''5.50'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.14%  '

$ws.Range("B17").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C17").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D17").Value = '3.852.76'
$ws.Range("E17").Value = '  +5.10%  '

$ws.Range("D18").Value = '3.299.15'
$ws.Range("E18").Value = '  +6.38%  '

$ws.Range("E19").Value = '  -0.70%  '

$ws.Range("D20").Value = '''498.61'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +12.94%  '

$ws.Range("D21").Value = '''15.15'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.52%  '

$ws.Range("D22").Value = '''5.93'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.53%  '

$ws.Range("D23").Value = '''0.0000209'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +5.70%  '

$ws.Range("D24").Value = '''9.18'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.50%  '

$ws.Range("D25").Value = '''5.62'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.03%  '

$ws.Range("D26").Value = '''88.87'
$ws.Range("D26").ClearFormats()

$ws.Range("D27").Value = '''12.01'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.61%  '

$ws.Range("D28").Value = '3.414.88'
$ws.Range("E28").Value = '  +4.85%  '

$ws.Range("E29").Value = '  +0.10%  '

$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").Value = '''0.245'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.67%  '

$ws.Range("B31").Value = 'Cronos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D31").Value = '''0.181'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.50%  '

$ws.Range("D32").Value = '''0.123'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.32%  '

$ws.Range("E33").Value = '  -0.08%  '

$ws.Range("D34").Value = '''9.30'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +2.01%  '

$ws.Range("D35").Value = '''27.61'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +7.20%  '

$ws.Range("E36").Value = '  -0.85%  '

$ws.Range("D37").Value = '''7.48'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -5.35%  '

$ws.Range("D38").Value = '''1.93'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +2.58%  '

$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").Value = '''492.15'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +3.77%  '

$ws.Range("B40").Value = 'WhiteBITCoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D40").Value = '''24.56'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +2.43%  '

$ws.Range("D41").Value = '''0.448'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.28%  '

$ws.Range("D42").Value = '''1.26'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.48%  '

$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").Value = '''3.27'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +2.16%  '

$ws.Range("B44").Value = 'MantraDAO'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D44").Value = '''3.50'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -9.26%  '

$ws.Range("E45").Value = '  +0.01%  '

$ws.Range("B46").Value = 'ARBITRUM'
$ws.Range("C46").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D46").Value = '''0.739'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +8.16%  '

$ws.Range("B47").Value = 'Monero'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D47").Value = '''161.56'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.03%  '

$ws.Range("D48").Value = '''1.95'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +5.94%  '

$ws.Range("D49").Value = '''4.60'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +6.01%  '

$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").Value = '''0.0336'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +5.36%  '

$ws.Range("D51").Value = '''0.792'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +9.22%  '
